$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.097.85'
$ws.Range("E2").Value = '  -0.71%  '

# Row 3
$ws.Range("D3").Value = '1.652.32'
$ws.Range("E3").Value = '  -0.75%  '

# Row 4
$ws.Range("E4").Value = '  -0.43%  '

# Row 5
$ws.Range("D5").Value = '''218.10'
$ws.Range("E5").Value = '  -0.29%  '

# Row 6
$ws.Range("D6").Value = '''0.5204'
$ws.Range("E6").Value = '  -2.22%  '

# Row 8
$ws.Range("D8").Value = '''0.2621'
$ws.Range("E8").Value = '  -1.02%  '

# Row 9
$ws.Range("D9").Value = '''0.06294'
$ws.Range("E9").Value = '  -1.36%  '

# Row 10
$ws.Range("D10").Value = '''20.50'
$ws.Range("E10").Value = '  -0.23%  '

# Row 11
$ws.Range("D11").Value = '''0.07812'
$ws.Range("E11").Value = '  -0.12%  '

# Row 12
$ws.Range("D12").Value = '''4.476'
$ws.Range("E12").Value = '  -1.69%  '

# Row 13
$ws.Range("D13").Value = '1.657.68'
$ws.Range("E13").Value = '  -0.35%  '

# Row 14
$ws.Range("D14").Value = '1.878.75'
$ws.Range("E14").Value = '  -0.74%  '

# Row 15
$ws.Range("D15").Value = '''0.5545'
$ws.Range("E15").Value = '  +0.48%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8015'
$ws.Range("E16").Value = '  -2.37%  '

# Row 17
$ws.Range("D17").Value = '''64.88'
$ws.Range("E17").Value = '  -1.15%  '

# Row 18
$ws.Range("D18").Value = '26.093.07'
$ws.Range("E18").Value = '  -0.85%  '

# Row 19
$ws.Range("E19").Value = '  -0.45%  '

# Row 20
$ws.Range("D20").Value = '''4.631'
$ws.Range("E20").Value = '  -1.17%  '

# Row 21
$ws.Range("D21").Value = '''194.62'
$ws.Range("E21").Value = '  +0.58%  '

# Row 22
$ws.Range("D22").Value = '''10.08'
$ws.Range("E22").Value = '  -1.24%  '

# Row 23
$ws.Range("D23").Value = '''5.952'
$ws.Range("E23").Value = '  -1.29%  '

# Row 24
$ws.Range("D24").Value = '''1.006'
$ws.Range("E24").Value = '  -0.44%  '

# Row 25
$ws.Range("D25").Value = '''146.61'
$ws.Range("E25").Value = '  +0.68%  '

# Row 26
$ws.Range("D26").Value = '''0.1205'
$ws.Range("E26").Value = '  -1.92%  '

# Row 27
$ws.Range("D27").Value = '''7.171'
$ws.Range("E27").Value = '  -0.39%  '

# Row 28
$ws.Range("D28").Value = '''15.90'
$ws.Range("E28").Value = '  -1.34%  '

# Row 29
$ws.Range("D29").Value = '''1.478'
$ws.Range("E29").Value = '  -0.14%  '

# Row 30
$ws.Range("D30").Value = '''0.05694'
$ws.Range("E30").Value = '  -2.96%  '

# Row 31
$ws.Range("E31").Value = '  -1.18%  '

# Row 32
$ws.Range("D32").Value = '''3.486'
$ws.Range("E32").Value = '  -3.48%  '

# Row 33
$ws.Range("D33").Value = '''3.384'
$ws.Range("E33").Value = '  +3.23%  '

# Row 34
$ws.Range("E34").Value = '  -1.13%  '

# Row 35
$ws.Range("D35").Value = '''2.799'
$ws.Range("E35").Value = '  -0.90%  '

# Row 36
$ws.Range("D36").Value = '''0.9508'
$ws.Range("E36").Value = '  -1.28%  '

# Row 37
$ws.Range("D37").Value = '''2.410'
$ws.Range("E37").Value = '  -0.34%  '

# Row 38
$ws.Range("D38").Value = '''0.5663'
$ws.Range("E38").Value = '  -2.39%  '

# Row 39
$ws.Range("D39").Value = '''0.01589'
$ws.Range("E39").Value = '  -1.24%  '

# Row 40
$ws.Range("D40").Value = '''5.958'
$ws.Range("E40").Value = '  +2.19%  '

# Row 41
$ws.Range("D41").Value = '1.056.90'
$ws.Range("E41").Value = '  +0.57%  '

# Row 42
$ws.Range("D42").Value = '''1.005'
$ws.Range("E42").Value = '  -0.39%  '

# Row 43
$ws.Range("D43").Value = '''0.8415'
$ws.Range("E43").Value = '  -2.83%  '

# Row 44
$ws.Range("D44").Value = '''103.69'
$ws.Range("E44").Value = '  -0.73%  '

# Row 45
$ws.Range("D45").Value = '1.790.18'
$ws.Range("E45").Value = '  -0.73%  '

# Row 46
$ws.Range("D46").Value = '''57.32'

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₈106'
$ws.Range("E47").Value = '  +1.60%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.05388'
$ws.Range("E48").Value = '  +4.36%  '

# Row 49
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = '''1.006'
$ws.Range("E49").Value = '  +0.10%  '

# Row 50
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.4393'
$ws.Range("E50").Value = '  +0.25%  '

# Row 51
$ws.Range("D51").Value = '''7.956'
$ws.Range("E51").Value = '  -1.31%  '
